$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row for 09/15/2025 (appended after the last existing row, 13)
$dateCell = $ws.Range("A14")
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/15/2025"
$dateCell.ClearFormats()

$ws.Range("B14").Value = 0.1268672426830251
$ws.Range("C14").Value = 0.8731327573169749
